# Weekly data refresh: insert a new weekly observation at row 134,
# shifting all subsequent rows (134-237) down by one (to 135-238).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 134 - this pushes the existing
# rows 134..237 down to 135..238 and grows the sheet dimension to R238.
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new weekly record.
$ws.Range("A134").Value = 9
$ws.Range("B134").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C134").Value = "Metropolitana"
$ws.Range("D134").Value = 44606
$ws.Range("E134").Value = 13
$ws.Range("F134").Value = 100112043
$ws.Range("G134").Value = "Pepino ensalada"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 79
$ws.Range("K134").Value = 14000
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = 14494
$ws.Range("N134").Value = "$/caja 60 unidades"
$ws.Range("O134").Value = "Región de Arica y Parinacota"
$ws.Range("P134").Value = 242
$ws.Range("Q134").Value = 60
$ws.Range("R134").Value = "Hortaliza"
